# Update the "dSF" (column F) values for several rows, per the commit:
# "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F6").Value  = -1
$ws.Range("F16").Value = 2
$ws.Range("F17").Value = 0
$ws.Range("F18").Value = -5
$ws.Range("F27").Value = -4
